# [Fonds de solidarite] Add 2021-01-22 data
#
# This script updates the "nombre_aides" (C), "nombre_entreprises" (D) and
# "montant_total" (E) columns for a handful of rows in the single sheet of
# the workbook, reflecting newly added/updated source data for 2021-01-22.
#
# All cells in this sheet are stored as text (inline strings) rather than
# numbers, so we force the NumberFormat to Text ("@") before writing the
# new values - this keeps Excel from silently re-typing them as numeric
# cells, matching the original data layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Worksheet, $Address, $Value)
    $rng = $Worksheet.Range($Address)
    $rng.NumberFormat = "@"
    $rng.Value = $Value
}

# Row 4
Set-TextValue $ws "C4" "671"
Set-TextValue $ws "D4" "566"
Set-TextValue $ws "E4" "6908275.49"

# Row 25
Set-TextValue $ws "C25" "85"
Set-TextValue $ws "E25" "1379797.91"

# Row 37
Set-TextValue $ws "C37" "386"
Set-TextValue $ws "E37" "3350933.69"

# Row 62
Set-TextValue $ws "C62" "1140"
Set-TextValue $ws "E62" "3803659.28"

# Row 63
Set-TextValue $ws "C63" "5716"
Set-TextValue $ws "E63" "24531453.78"

# Row 64
Set-TextValue $ws "C64" "3149"
Set-TextValue $ws "E64" "19489380.84"

# Row 65
Set-TextValue $ws "C65" "1116"
Set-TextValue $ws "D65" "1012"
Set-TextValue $ws "E65" "9135015.17"

# Row 66
Set-TextValue $ws "C66" "317"
Set-TextValue $ws "D66" "286"
Set-TextValue $ws "E66" "4380004.41"

# Row 67
Set-TextValue $ws "C67" "51"
Set-TextValue $ws "E67" "1887881.18"

# Row 94
Set-TextValue $ws "C94" "488"
Set-TextValue $ws "E94" "4743869.46"
